# Update Grn-Tnfrsf1b.xlsx LR-pair stats with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns G:J (Ligand avg/total expression + derived specificity of avg/total)
$valsGJ = @(
    @(33.36960033333333, 100.108801, 0.07727383968381614, 0.07727383968381614),
    @(33.36960033333333, 100.108801, 0.07727383968381614, 0.07727383968381614),
    @(33.36960033333333, 100.108801, 0.07727383968381614, 0.07727383968381614),
    @(33.36960033333333, 100.108801, 0.07727383968381614, 0.07727383968381614),
    @(110.1980973333333, 330.594292, 0.2551852590901843, 0.2551852590901843),
    @(110.1980973333333, 330.594292, 0.2551852590901843, 0.2551852590901843),
    @(110.1980973333333, 330.594292, 0.2551852590901843, 0.2551852590901843),
    @(110.1980973333333, 330.594292, 0.2551852590901843, 0.2551852590901843),
    @(13.90116633333333, 41.703499, 0.0321908709702775, 0.0321908709702775),
    @(13.90116633333333, 41.703499, 0.0321908709702775, 0.0321908709702775),
    @(13.90116633333333, 41.703499, 0.0321908709702775, 0.0321908709702775),
    @(13.90116633333333, 41.703499, 0.0321908709702775, 0.0321908709702775),
    @(274.366806, 823.100418, 0.635350030255722, 0.635350030255722),
    @(274.366806, 823.100418, 0.635350030255722, 0.635350030255722),
    @(274.366806, 823.100418, 0.635350030255722, 0.635350030255722),
    @(274.366806, 823.100418, 0.635350030255722, 0.635350030255722)
)

# Columns M:T (Receptor avg/total expr + specificity, edge avg/total weight + specificity)
$valsMT = @(
    @(12.57753066666667, 37.732592, 0.1317204310459389, 0.1317204310459389, 419.7071715269102, 3777.364543742191, 0.01017854347172704, 0.01017854347172704),
    @(10.940628, 32.821884, 0.1145776761962127, 0.1145776761962127, 365.0843837556759, 3285.759453801084, 0.008853856981730334, 0.008853856981730334),
    @(2.520244333333333, 7.560733, 0.02639370785296846, 0.02639370785296846, 84.09954614568144, 756.895915311133, 0.002039543149291764, 0.002039543149291764),
    @(69.44815566666666, 208.344467, 0.7273081849048799, 0.7273081849048799, 2317.457198483785, 20857.11478635407, 0.05620189608106699, 0.05620189608106699),
    @(12.57753066666667, 37.732592, 0.1317204310459389, 0.1317204310459389, 1386.019948618318, 12474.17953756486, 0.03361311232392867, 0.03361311232392867),
    @(10.940628, 32.821884, 0.1145776761962127, 0.1145776761962127, 1205.636389231792, 10850.72750308613, 0.02923853398608177, 0.02923853398608177),
    @(2.520244333333333, 7.560733, 0.02639370785296846, 0.02639370785296846, 277.7261303484485, 2499.535173136036, 0.006735285176810388, 0.006735285176810388),
    @(69.44815566666666, 208.344467, 0.7273081849048799, 0.7273081849048799, 7653.054617775818, 68877.49155998237, 0.1855983276033634, 0.1855983276033634),
    @(12.57753066666667, 37.732592, 0.1317204310459389, 0.1317204310459389, 174.8423458599342, 1573.581112739408, 0.004240195399949153, 0.004240195399949153),
    @(10.940628, 32.821884, 0.1145776761962127, 0.1145776761962127, 152.087489619124, 1368.787406572116, 0.003688355190506518, 0.003688355190506518),
    @(2.520244333333333, 7.560733, 0.02639370785296846, 0.02639370785296846, 35.03433567830745, 315.309021104767, 0.0008496364439221077, 0.0008496364439221077),
    @(69.44815566666666, 208.344467, 0.7273081849048799, 0.7273081849048799, 965.4103634655593, 8688.693271190034, 0.02341268393589972, 0.02341268393589972),
    @(12.57753066666667, 37.732592, 0.1317204310459389, 0.1317204310459389, 3450.856916380384, 31057.71224742345, 0.08368857985033402, 0.08368857985033402),
    @(10.940628, 32.821884, 0.1145776761962127, 0.1145776761962127, 3001.745159994167, 27015.70643994751, 0.07279693003789404, 0.07279693003789404),
    @(2.520244333333333, 7.560733, 0.02639370785296846, 0.02639370785296846, 691.471388076266, 6223.242492686394, 0.0167692430829442, 0.0167692430829442),
    @(69.44815566666666, 208.344467, 0.7273081849048799, 0.7273081849048799, 19054.26865285413, 171488.4178756872, 0.4620952772845497, 0.4620952772845497)
)

$colsGJ = @("G","H","I","J")
$colsMT = @("M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    for ($c = 0; $c -lt $colsGJ.Length; $c++) {
        $ws.Range("$($colsGJ[$c])$row").Value2 = $valsGJ[$i][$c]
    }
    for ($c = 0; $c -lt $colsMT.Length; $c++) {
        $ws.Range("$($colsMT[$c])$row").Value2 = $valsMT[$i][$c]
    }
}
